$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$values = @(
    "-4,-7",
    "-1,-8",
    "-9,-9",
    "-10,-2",
    "8,4",
    "-3,-6",
    "-3,3",
    "3,-6",
    "4,1",
    "1,-1",
    "5,-7",
    "6,5",
    "-3,2",
    "-2,-9",
    "0,8",
    "6,5",
    "2,4",
    "2,-8",
    "-10,-5",
    "-1,6"
)

$row = 2
foreach ($val in $values) {
    $ws.Range("D$row").Value = $val
    $row++
}
